# LOM3077.xlsx update: fill in the real "Usinagem de Materiais" course data.
# The sheet had its B/C (Portuguese/English) content columns out of sync with
# their A-column labels (looks like a row got dropped upstream). Fix it by
# inserting the missing "Docentes responsaveis" data row and then restoring
# the correct long-form text beside each label.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert the missing row under "Docentes responsaveis:" (row 12) ---
# This shifts everything from row 13 down by one, which is exactly what the
# target layout needs (row 13 becomes a B/C-only row holding the professor
# name, and every label below slides down to line up with its real text).
$ws.Rows.Item(13).Insert()

# The freshly inserted row 13 inherited row-12's formatting (a leftover A13
# cell). Clear it completely so it starts blank, matching the target sheet.
$ws.Range("A13:C13").Clear()

# --- 2. Professor name moves into the new row 13 (B13/C13 only) ---
# Borrow the wrap-text/vertical-top formatting already used by column B/C
# elsewhere on the sheet (e.g. row 10) instead of re-deriving it, so the
# cells land on the existing shared style rather than minting a new one.
$ws.Range("B10").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C10").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$profName = "471420 - Carlos Antonio Reis Pereira Baptista"
$ws.Range("B13").Value = $profName
$ws.Range("C13").Value = $profName

# --- 3. Objetivos / Objectives body text (row 10) ---
$objText = @"
Desenvolver a capacidade de elaborar, interpretar e executar processos de usinagem. Apresentar e discutir conceitos sobre os processos de usinagem e sobre as interações devido ao contato ferramenta-peça durante as operações de corte. Apresentar os mecanismos de desgaste e avaria além da correlação com a vida da ferramenta. Propor abordagens com enfoque nos aspectos econômicos da usinagem dos materiais. Discutir sobre as novas tendências da tecnologia da usinagem
"@
$ws.Range("B10").Value = $objText
$ws.Range("C10").Value = $objText

# --- 4. Programa resumido / Short syllabus body text (row 14) ---
$shortSyllabus = @"
Usinagem: tipos, propriedades, processamento e aplicações. Tendências da tecnologia da usinagem.
"@
$ws.Range("B14").Value = $shortSyllabus
$ws.Range("C14").Value = $shortSyllabus

# --- 5. Programa / Syllabus body text (row 16) ---
$programa = @"
Conteúdo teórico: 
1. Processos mecânicos de usinagem. 
2. Mecanismos de formação do cavaco. 
3. Cálculo de potência de corte. 
4. Condições econômicas de corte.
5. Fluido de corte.
6. Usinabilidade dos materiais.
7. Novas tendências na usinagem dos materiais. 
Conteúdo prático: 
1. Trabalhos práticos em máquinas e equipamentos de usinagem.
2. Visita a empresa de usinagem.
"@
$ws.Range("B16").Value = $programa
$ws.Range("C16").Value = $programa

# --- 5b. Método / Critério / Norma de recuperação body text (rows 19-21) ---
# These three rows were already mis-mapped one below where they belonged
# before this edit (pre-existing data-entry bug), so fix all three in place.
$metodoText = @"
Aula expositiva com utilização de recursos audiovisuais, aliada a aulas práticas de preparação de máquinas e equipamentos de usinagem e visita a empresas.
"@
$ws.Range("B19").Value = $metodoText
$ws.Range("C19").Value = $metodoText

$criterioText = @"
A nota final será calculada pela média ponderada de duas provas, valendo 60% e da média de exercícios, testes práticos e relatórios de laboratório, valendo 40% da nota final.
A fórmula para o cálculo da média será: NF = (P1+2*P2)/3*0,6 + ME*0,4, na qual P1 e P2 são as notas das provas e ME a média dos exercícios, testes e relatórios.
"@
$ws.Range("B20").Value = $criterioText
$ws.Range("C20").Value = $criterioText

$recuperacaoText = @"
A recuperação será uma prova dissertativa. A média final será MF = (NF + RE)/2, na qual NF é a nota final e RE a nota da prova de recuperação.
"@
$ws.Range("B21").Value = $recuperacaoText
$ws.Range("C21").Value = $recuperacaoText

# --- 6. Bibliografia body text (row 22) ---
$biblio = @"
1.SCHNEIDER, JR. G., Cutting Tool Applications. Nelson Publishing, Inc. New York, USA, 2001. 
2.DINIZ, A.E.; MARCONDES, F.C.; COPPINI, N.L., Tecnologia da Usinagem dos Materiais. Ed. Artliber, São Paulo, 2ª ed., 2000.
3.FERRAREZI, Dino. Fundamentos da usinagem dos metais, Edgar Blucher, 1995.
4.ABNT - Normas Técnicas de 1995 - edição ABNT.
5.SANDVIK COROMANT, Modern Metal Cutting. AB Sandvik Coromant, Sandviken, Sweden, 1994. 
6.MACHADO, A. Usinagem dos metais. Uberlândia: Universidade Federal de Uberlândia, 1994.
7.DeVRIES, W.R., Analysis of Material Removal Processes. Springer-Verlag, New York, USA, 1991.
"@
$ws.Range("B22").Value = $biblio
$ws.Range("C22").Value = $biblio

Write-Host "done"
